$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the existing "firstname" field label to "author.firstname"
$ws.Range("A5").Value = "author.firstname"

# Fix category.2 value from "Technology" to "Non-Fiction"
$ws.Range("C7").Value = "Non-Fiction"

# Add a new "Record-2" column (column D) with its values
$ws.Range("D1").Value = "Record-2"
$ws.Range("D2").Value = "234-456-2210"
$ws.Range("D3").Value = 12345678
$ws.Range("D4").Value = "Shawn"
$ws.Range("D5").Value = "Peter"
$ws.Range("D6").Value = "Non-Fiction"
$ws.Range("D7").Value = "Non-Fiction"
$ws.Range("D8").Value = 1000

# Update selection to match the new active cell
$ws.Range("D8").Select()
